# Apply data refresh to 杭州-漫展信息.xlsx
# Sheets (in workbook order): 1=展览, 2=演出, 3=本地生活, 4=全部类型
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# ---- Sheet 1 (展览) : column F ("想去人数") updates ----
$sheet1F = @{
  2  = 294
  4  = 9480
  5  = 202
  6  = 72
  7  = 1966
  8  = 6484
  9  = 628
  11 = 9937
  12 = 11409
  13 = 1246
  14 = 1173
  15 = 4983
  16 = 808
  17 = 474
  21 = 1348
  22 = 267
  23 = 1877
  24 = 903
  25 = 1281
  28 = 2062
  29 = 435
  30 = 645
  31 = 2702
  32 = 191
  33 = 1789
  34 = 96
  35 = 809
  36 = 68
  37 = 924
  39 = 3374
  41 = 89
  45 = 901
  46 = 247
  48 = 4224
  49 = 62
}
foreach ($row in $sheet1F.Keys) {
  $ws1.Cells.Item($row, 6).Value = $sheet1F[$row]
}

# Sheet 1 : column I ("Cover") updated cover image URL for row 35
$ws1.Cells.Item(35, 9).Value = "//i1.hdslb.com/bfs/openplatform/202407/kZIRM1Sx1720073676616.jpeg"

# ---- Sheet 2 (演出) : column F updates ----
$sheet2F = @{
  9  = 19
  25 = 106
}
foreach ($row in $sheet2F.Keys) {
  $ws2.Cells.Item($row, 6).Value = $sheet2F[$row]
}

# ---- Sheet 3 (本地生活) : column F updates ----
$ws3.Cells.Item(2, 6).Value = 6009

# ---- Sheet 4 (全部类型) : column F updates ----
$sheet4F = @{
  2  = 294
  4  = 9480
  5  = 202
  6  = 72
  10 = 6484
  11 = 9937
  12 = 11409
  13 = 1173
  14 = 808
  15 = 474
  18 = 19
  22 = 1348
  23 = 267
  24 = 1877
  25 = 903
  26 = 1281
  28 = 2062
  29 = 645
  30 = 2702
  31 = 191
  32 = 96
  33 = 809
  38 = 68
  42 = 89
  45 = 901
  46 = 247
  48 = 4224
}
foreach ($row in $sheet4F.Keys) {
  $ws4.Cells.Item($row, 6).Value = $sheet4F[$row]
}

# Sheet 4 : column I ("Cover") updated cover image URL for row 33
$ws4.Cells.Item(33, 9).Value = "//i1.hdslb.com/bfs/openplatform/202407/kZIRM1Sx1720073676616.jpeg"
